$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '63.698.11'
$ws.Range("E2").Value2 = '  -0.83%  '
$ws.Range("D3").Value2 = '3.120.62'
$ws.Range("E3").Value2 = '  -1.23%  '
$ws.Range("E4").Value2 = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.01'
$ws.Range("E5").Value2 = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.80'
$ws.Range("E6").Value2 = '  -3.41%  '
$ws.Range("E7").Value2 = '  -0.10%  '
$ws.Range("D8").Value2 = '3.120.20'
$ws.Range("E8").Value2 = '  -1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value2 = '  -0.85%  '
$ws.Range("E10").Value2 = '  -2.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.37'
$ws.Range("E11").Value2 = '  -0.97%  '
$ws.Range("E12").Value2 = '  -1.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value2 = '  -2.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.15'
$ws.Range("E14").Value2 = '  -1.24%  '
$ws.Range("D15").Value2 = '3.630.41'
$ws.Range("E15").Value2 = '  -1.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.119'
$ws.Range("E16").Value2 = '  +2.30%  '
$ws.Range("D17").Value2 = '63.706.97'
$ws.Range("E17").Value2 = '  -0.78%  '
$ws.Range("D18").Value2 = '3.115.74'
$ws.Range("E18").Value2 = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.83'
$ws.Range("E19").Value2 = '  -1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.98'
$ws.Range("E20").Value2 = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.78'
$ws.Range("E21").Value2 = '  +0.26%  '
$ws.Range("E22").Value2 = '  -1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.61'
$ws.Range("E23").Value2 = '  -5.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.60'
$ws.Range("E24").Value2 = '  +3.33%  '
$ws.Range("E25").Value2 = '  -3.02%  '
$ws.Range("E26").Value2 = '  +0.03%  '
$ws.Range("E27").Value2 = '  -3.34%  '
$ws.Range("E28").Value2 = '  -3.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.94'
$ws.Range("E29").Value2 = '  -3.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("E30").Value2 = '  -3.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.08'
$ws.Range("E31").Value2 = '  +2.61%  '
$ws.Range("E32").Value2 = '  -0.10%  '
$ws.Range("E33").Value2 = '  -8.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.66'
$ws.Range("E34").Value2 = '  -2.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value2 = '  -2.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("E36").Value2 = '  -0.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.64'
$ws.Range("E37").Value2 = '  -0.67%  '
$ws.Range("D38").Value2 = '0.0₃0742'
$ws.Range("E38").Value2 = '  -5.96%  '
$ws.Range("E39").Value2 = '  -7.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '436.50'
$ws.Range("E40").Value2 = '  -5.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0394'
$ws.Range("E41").Value2 = '  -1.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.120'
$ws.Range("E42").Value2 = '  +0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.28'
$ws.Range("E43").Value2 = '  -1.27%  '
$ws.Range("D44").Value2 = '2.867.08'
$ws.Range("E44").Value2 = '  -0.12%  '
$ws.Range("E45").Value2 = '  -3.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("E46").Value2 = '  -4.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.39'
$ws.Range("E47").Value2 = '  -2.00%  '
$ws.Range("E48").Value2 = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.87'
$ws.Range("E49").Value2 = '  -2.83%  '
$ws.Range("E50").Value2 = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.58'
$ws.Range("E51").Value2 = '  +1.97%  '
